# Update cryptos list: apply latest price/volume(1h) snapshot, and reorder a few coin rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the data range (prices/percentages must stay text, not be
# reinterpreted as numbers) while setting the new values, then restore the default style
# so the cells end up with no explicit style, matching the original formatting.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = '66.460.47'
$ws.Cells.Item(2, 5).Value = '  -5.55%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.352.62'
$ws.Cells.Item(3, 5).Value = '  -6.71%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.09%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '558.08'
$ws.Cells.Item(5, 5).Value = '  -5.99%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '182.24'
$ws.Cells.Item(6, 5).Value = '  -8.97%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '0.596'
$ws.Cells.Item(7, 5).Value = '  -5.13%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.10%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '3.339.90'
$ws.Cells.Item(9, 5).Value = '  -6.77%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '0.187'
$ws.Cells.Item(10, 5).Value = '  -13.70%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '0.591'
$ws.Cells.Item(11, 5).Value = '  -8.26%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '47.69'
$ws.Cells.Item(12, 5).Value = '  -10.64%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '0.0000268'
$ws.Cells.Item(13, 5).Value = '  -10.88%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '8.68'
$ws.Cells.Item(14, 5).Value = '  -9.83%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '3.890.60'
$ws.Cells.Item(15, 5).Value = '  -6.67%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '596.73'
$ws.Cells.Item(16, 5).Value = '  -14.70%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '66.204.43'
$ws.Cells.Item(17, 5).Value = '  -6.05%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '3.361.94'
$ws.Cells.Item(18, 5).Value = '  -6.47%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '0.117'
$ws.Cells.Item(19, 5).Value = '  -4.51%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '17.69'
$ws.Cells.Item(20, 5).Value = '  -7.15%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '11.62'
$ws.Cells.Item(21, 5).Value = '  -8.80%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '0.908'
$ws.Cells.Item(22, 5).Value = '  -8.60%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '16.72'
$ws.Cells.Item(23, 5).Value = '  -8.13%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '5.02'
$ws.Cells.Item(24, 5).Value = '  -5.56%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '96.59'
$ws.Cells.Item(25, 5).Value = '  -13.03%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '4.06'
$ws.Cells.Item(26, 5).Value = '  -10.43%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '2.74'
$ws.Cells.Item(27, 5).Value = '  -8.64%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '9.44'
$ws.Cells.Item(28, 5).Value = '  -10.02%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '8.76'
$ws.Cells.Item(29, 5).Value = '  -12.26%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '30.62'
$ws.Cells.Item(30, 5).Value = '  -11.66%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '3.87'
$ws.Cells.Item(31, 5).Value = '  -12.52%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '6.31'
$ws.Cells.Item(32, 5).Value = '  -10.65%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '11.13'
$ws.Cells.Item(33, 5).Value = '  -9.09%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '0.105'
$ws.Cells.Item(34, 5).Value = '  -7.57%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'Maker'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(35, 4).Value = '3.826.89'
$ws.Cells.Item(35, 5).Value = '  +0.66%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'OKB'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(36, 4).Value = '58.21'
$ws.Cells.Item(36, 5).Value = '  -8.43%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'Bittensor'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(37, 4).Value = '531.31'
$ws.Cells.Item(37, 5).Value = '  +4.27%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '0.999'
$ws.Cells.Item(38, 5).Value = '  -0.04%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '3.69'
$ws.Cells.Item(39, 5).Value = '  +34.64%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '3.38'
$ws.Cells.Item(40, 5).Value = '  -7.01%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '0.0₃0722'
$ws.Cells.Item(41, 5).Value = '  -14.69%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'Kaspa'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(42, 4).Value = '0.128'
$ws.Cells.Item(42, 5).Value = '  -6.63%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'Fetch.AI'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(43, 4).Value = '2.69'
$ws.Cells.Item(43, 5).Value = '  -10.35%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '0.348'
$ws.Cells.Item(44, 5).Value = '  -8.67%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '32.39'
$ws.Cells.Item(45, 5).Value = '  -11.08%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '0.0416'
$ws.Cells.Item(46, 5).Value = '  -11.57%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '2.66'
$ws.Cells.Item(47, 5).Value = '  -12.53%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '3.12'
$ws.Cells.Item(48, 5).Value = '  -9.83%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '0.130'
$ws.Cells.Item(49, 5).Value = '  -8.10%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '1.00'
$ws.Cells.Item(50, 5).Value = '  -0.16%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '7.69'
$ws.Cells.Item(51, 5).Value = '  -11.23%  '

$dataRange.Style = "Normal"
